# Apply Harvard case classification results to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (precision)
$ws.Range("C2").Value = 0.4
$ws.Range("M2").Value = 0.25

# Row 3 (recall)
$ws.Range("C3").Value = 0.5
$ws.Range("M3").Value = 0.25

# Row 4 (f1-score)
$ws.Range("C4").Value = 0.4444444444444445
$ws.Range("M4").Value = 0.25

# Row 5 (f2-score)
$ws.Range("C5").Value = 0.4761904761904762
$ws.Range("M5").Value = 0.25

# Row 6 (NDCG)
$ws.Range("C6").Value = 0.3604963103757638
$ws.Range("M6").Value = 0.07026687648649678
